$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.088543891906738
$ws.Range("B1").Value = 2.970503568649292
$ws.Range("C1").Value = 2.688413619995117
$ws.Range("D1").Value = 3.012854099273682
$ws.Range("E1").Value = 2.758003234863281
